# Update "想去人数" (number of people interested) values to match the
# newly generated output at commit 456a3b4.
#
# Affected sheets: 展览 (Exhibitions) and 全部类型 (All types), which both
# list the same events (演出/本地生活 are unaffected by this data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1182   # 南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展: 1180 -> 1182
$wsExhibit.Range("F5").Value = 4992   # 南宁·AB动漫游戏嘉年华: 0 -> 4992
$wsExhibit.Range("F7").Value = 9302   # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）: 9298 -> 9302
$wsExhibit.Range("F11").Value = 649   # 南宁·第二届北极光动漫展: 645 -> 649

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1182    # 南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展: 1180 -> 1182
$wsAll.Range("F8").Value = 523     # 横州·第二届海棠动漫游戏嘉年华: 0 -> 523
$wsAll.Range("F10").Value = 9302   # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）: 9298 -> 9302
$wsAll.Range("F16").Value = 649    # 南宁·第二届北极光动漫展: 645 -> 649

$wb.Save()
